# Auto-generated edit script: updates market-price snapshot values (H..N columns)
# across multiple Leve sheets, reproducing the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (42 cell changes) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 896.13336
$ws.Range("I9").Value = 966.55554
$ws.Range("J9").Value = 790.5
$ws.Range("K9").Value = 966.55554
$ws.Range("L9").Value = 790.5
$ws.Range("M9").Value = -797.55554
$ws.Range("N9").Value = -1128.5
$ws.Range("H12").Value = 698
$ws.Range("I12").Value = 698
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 698
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -528
$ws.Range("N12").ClearContents()
$ws.Range("H62").Value = 3674.75
$ws.Range("I62").Value = 3299.6667
$ws.Range("K62").Value = 3299.6667
$ws.Range("M62").Value = -2675.6667
$ws.Range("H65").Value = 3674.75
$ws.Range("I65").Value = 3299.6667
$ws.Range("K65").Value = 16498.3335
$ws.Range("M65").Value = -13378.3335
$ws.Range("H115").Value = 380
$ws.Range("I115").Value = 380
$ws.Range("K115").Value = 1140
$ws.Range("M115").Value = 427
$ws.Range("H135").Value = 453.875
$ws.Range("I135").Value = 453.875
$ws.Range("K135").Value = 4084.875
$ws.Range("M135").Value = -1549.875
$ws.Range("H137").Value = 1182.7
$ws.Range("I137").Value = 1286.5333
$ws.Range("K137").Value = 3859.5999
$ws.Range("M137").Value = -1309.5999
$ws.Range("H138").Value = 858.5
$ws.Range("I138").Value = 858.5
$ws.Range("K138").Value = 2575.5
$ws.Range("M138").Value = 2564.5
$ws.Range("H141").Value = 9986.25
$ws.Range("I141").Value = 9986.25
$ws.Range("K141").Value = 29958.75
$ws.Range("M141").Value = -24778.75

# --- Sheet: ARM (18 cell changes) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 590.8421
$ws.Range("I74").Value = 548.64703
$ws.Range("J74").Value = 949.5
$ws.Range("K74").Value = 548.64703
$ws.Range("L74").Value = 949.5
$ws.Range("M74").Value = 325.35297
$ws.Range("N74").Value = -2697.5
$ws.Range("H77").Value = 590.8421
$ws.Range("I77").Value = 548.64703
$ws.Range("J77").Value = 949.5
$ws.Range("K77").Value = 2743.23515
$ws.Range("L77").Value = 4747.5
$ws.Range("M77").Value = 1624.76485
$ws.Range("N77").Value = -13483.5
$ws.Range("H110").Value = 2855.8572
$ws.Range("I110").Value = 2855.8572
$ws.Range("K110").Value = 2855.8572
$ws.Range("M110").Value = -810.8571999999999

# --- Sheet: CRP (51 cell changes) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H3").Value = 221.5
$ws.Range("I3").Value = 221.5
$ws.Range("K3").Value = 221.5
$ws.Range("M3").Value = -108.5
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 508.7143
$ws.Range("I5").Value = 317
$ws.Range("J5").Value = 652.5
$ws.Range("K5").Value = 317
$ws.Range("L5").Value = 652.5
$ws.Range("M5").Value = -205
$ws.Range("N5").Value = -876.5
$ws.Range("H11").Value = 10000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10280
$ws.Range("H12").Value = 3826.3333
$ws.Range("I12").Value = 3826.3333
$ws.Range("K12").Value = 3826.3333
$ws.Range("M12").Value = -3656.3333
$ws.Range("H42").Value = 4499.5
$ws.Range("I42").Value = 4333.3335
$ws.Range("K42").Value = 4333.3335
$ws.Range("M42").Value = -3740.3335
$ws.Range("H51").Value = 3199.5
$ws.Range("I51").Value = 3199.5
$ws.Range("K51").Value = 3199.5
$ws.Range("M51").Value = -2463.5
$ws.Range("H59").Value = 85020.39999999999
$ws.Range("J59").Value = 147499
$ws.Range("L59").Value = 147499
$ws.Range("N59").Value = -149789
$ws.Range("H61").Value = 3199.5
$ws.Range("I61").Value = 3199.5
$ws.Range("K61").Value = 3199.5
$ws.Range("M61").Value = -2851.5
$ws.Range("H94").Value = 1049
$ws.Range("J94").Value = 573.5
$ws.Range("L94").Value = 573.5
$ws.Range("N94").Value = -1475.5
$ws.Range("H140").Value = 39999
$ws.Range("J140").Value = 39999
$ws.Range("L140").Value = 39999
$ws.Range("N140").Value = -50359

# --- Sheet: CUL (62 cell changes) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 84.90000000000001
$ws.Range("I2").Value = 74.454544
$ws.Range("J2").Value = 97.666664
$ws.Range("K2").Value = 446.727264
$ws.Range("L2").Value = 585.999984
$ws.Range("M2").Value = -333.727264
$ws.Range("N2").Value = -811.999984
$ws.Range("H6").Value = 157.44444
$ws.Range("I6").Value = 145.875
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 437.625
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -324.625
$ws.Range("N6").Value = -976
$ws.Range("H7").Value = 14950.071
$ws.Range("I7").Value = 17275.084
$ws.Range("K7").Value = 51825.25199999999
$ws.Range("M7").Value = -51713.25199999999
$ws.Range("H9").Value = 1186.5217
$ws.Range("I9").Value = 250
$ws.Range("J9").Value = 1275.7142
$ws.Range("K9").Value = 750
$ws.Range("L9").Value = 3827.1426
$ws.Range("M9").Value = -526
$ws.Range("N9").Value = -4275.142599999999
$ws.Range("H10").Value = 296.06668
$ws.Range("I10").Value = 203.28572
$ws.Range("J10").Value = 1595
$ws.Range("K10").Value = 609.85716
$ws.Range("L10").Value = 4785
$ws.Range("M10").Value = -470.85716
$ws.Range("N10").Value = -5063
$ws.Range("H11").Value = 522.5
$ws.Range("I11").Value = 522.5
$ws.Range("K11").Value = 1567.5
$ws.Range("M11").Value = -1427.5
$ws.Range("H13").Value = 154.2
$ws.Range("I13").Value = 30.25
$ws.Range("J13").Value = 650
$ws.Range("K13").Value = 90.75
$ws.Range("L13").Value = 1950
$ws.Range("M13").Value = 77.25
$ws.Range("N13").Value = -2286
$ws.Range("H15").Value = 384.7
$ws.Range("J15").Value = 594.6667
$ws.Range("L15").Value = 1784.0001
$ws.Range("N15").Value = -2064.0001
$ws.Range("H17").Value = 235.2
$ws.Range("I17").Value = 77.42856999999999
$ws.Range("K17").Value = 232.28571
$ws.Range("M17").Value = -63.28570999999999
$ws.Range("H58").Value = 4000
$ws.Range("I58").Value = 4000
$ws.Range("K58").Value = 12000
$ws.Range("M58").Value = -11872
$ws.Range("H113").Value = 1435.1111
$ws.Range("I113").Value = 1398
$ws.Range("J113").Value = 1449.3846
$ws.Range("K113").Value = 4194
$ws.Range("L113").Value = 4348.1538
$ws.Range("M113").Value = -2024
$ws.Range("N113").Value = -8688.1538

# --- Sheet: GSM (7 cell changes) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 24067.5
$ws.Range("I24").Value = 20005
$ws.Range("J24").Value = 29290.715
$ws.Range("K24").Value = 20005
$ws.Range("L24").Value = 29290.715
$ws.Range("M24").Value = -19832
$ws.Range("N24").Value = -29636.715

# --- Sheet: LTW (11 cell changes) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1323.7142
$ws.Range("I40").Value = 1349.3846
$ws.Range("J40").Value = 990
$ws.Range("K40").Value = 1349.3846
$ws.Range("L40").Value = 990
$ws.Range("M40").Value = -1213.3846
$ws.Range("N40").Value = -1262
$ws.Range("H136").Value = 1110.4286
$ws.Range("I136").Value = 1110.4286
$ws.Range("K136").Value = 3331.2858
$ws.Range("M136").Value = -781.2857999999997

# --- Sheet: WVR (8 cell changes) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 3500
$ws.Range("I54").Value = 3500
$ws.Range("K54").Value = 3500
$ws.Range("M54").Value = -2980
$ws.Range("H136").Value = 4236.59
$ws.Range("I136").Value = 3917.35
$ws.Range("K136").Value = 11752.05
$ws.Range("M136").Value = -9202.049999999999
